$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (logistic_embeddings)
$ws.Range("C5").Value = 0.256
$ws.Range("D5").Value = 0.392
$ws.Range("E5").Value = 0.439
$ws.Range("F5").Value = 0.477
$ws.Range("G5").Value = 0.514
$ws.Range("H5").Value = 0.529

# Row 7 (classical-best-embeddings -> classical-best-embed)
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.256
$ws.Range("D7").Value = 0.392
$ws.Range("E7").Value = 0.439

# Row 8 (BERT-base)
$ws.Range("C8").Value = 0.246
$ws.Range("D8").Value = 0.48
$ws.Range("E8").Value = 0.522
$ws.Range("F8").Value = 0.583
$ws.Range("G8").Value = 0.625
$ws.Range("H8").Value = 0.633

# Row 9 (BERT-base-nli)
$ws.Range("B9").Value = 0.359
$ws.Range("C9").Value = 0.47
$ws.Range("D9").Value = 0.554
$ws.Range("E9").Value = 0.584
$ws.Range("G9").Value = 0.642
$ws.Range("H9").Value = 0.653
